$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '59.375.73'
Set-TextValue 'E2' '  +0.53%  '
Set-TextValue 'D3' '2.521.44'
Set-TextValue 'E3' '  +3.06%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '543.78'
Set-TextValue 'E5' '  +0.62%  '
Set-TextValue 'D6' '144.89'
Set-TextValue 'E6' '  -1.58%  '
Set-TextValue 'E7' '  -0.19%  '
Set-TextValue 'D8' '0.576'
Set-TextValue 'E8' '  +0.73%  '
Set-TextValue 'D9' '2.553.92'
Set-TextValue 'E9' '  +3.69%  '
Set-TextValue 'E10' '  +2.02%  '
Set-TextValue 'E11' '  +0.45%  '
Set-TextValue 'E12' '  +4.48%  '
Set-TextValue 'E13' '  +1.29%  '
Set-TextValue 'D14' '2.968.25'
Set-TextValue 'E14' '  +2.90%  '
Set-TextValue 'D15' '23.90'
Set-TextValue 'E15' '  -0.19%  '
Set-TextValue 'D16' '59.309.42'
Set-TextValue 'E16' '  +0.63%  '
Set-TextValue 'E17' '  +2.49%  '
Set-TextValue 'D18' '2.541.41'
Set-TextValue 'E18' '  +1.09%  '
Set-TextValue 'D19' '11.29'
Set-TextValue 'E19' '  +1.72%  '
Set-TextValue 'D20' '4.32'
Set-TextValue 'E20' '  -0.93%  '
Set-TextValue 'D21' '327.22'
Set-TextValue 'E21' '  +1.06%  '
Set-TextValue 'E22' '  +3.28%  '
Set-TextValue 'D23' '5.87'
Set-TextValue 'E23' '  +2.77%  '
Set-TextValue 'D24' '62.07'
Set-TextValue 'E24' '  +2.34%  '
Set-TextValue 'E25' '  -2.42%  '
Set-TextValue 'E26' '  +2.21%  '
Set-TextValue 'E27' '  +1.70%  '
Set-TextValue 'D28' '8.02'
Set-TextValue 'E28' '  +4.26%  '
Set-TextValue 'D29' '6.94'
Set-TextValue 'E29' '  +4.38%  '
Set-TextValue 'D30' '0.0₃0787'
Set-TextValue 'E30' '  +2.31%  '
Set-TextValue 'E31' '  +0.66%  '
Set-TextValue 'D32' '1.23'
Set-TextValue 'E32' '  -2.79%  '
Set-TextValue 'E33' '  +9.13%  '
Set-TextValue 'E34' '  -0.09%  '
Set-TextValue 'D35' '157.27'
Set-TextValue 'E35' '  +0.37%  '
Set-TextValue 'D36' '18.74'
Set-TextValue 'E36' '  +1.72%  '
Set-TextValue 'D37' '4.42'
Set-TextValue 'E37' '  -0.89%  '
Set-TextValue 'E38' '  -5.04%  '
Set-TextValue 'D39' '5.68'
Set-TextValue 'E39' '  -2.89%  '
Set-TextValue 'D40' '36.96'
Set-TextValue 'E40' '  +2.03%  '
Set-TextValue 'D41' '299.88'
Set-TextValue 'E41' '  -5.26%  '
Set-TextValue 'E42' '  +0.44%  '
Set-TextValue 'E43' '  -0.60%  '
Set-TextValue 'E44' '  -0.36%  '
Set-TextValue 'E45' '  +4.26%  '
Set-TextValue 'D46' '10.77'
Set-TextValue 'E46' '  +0.51%  '
Set-TextValue 'D47' '0.0936'
Set-TextValue 'E47' '  -0.62%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '18.87'
Set-TextValue 'E48' '  +2.58%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D49' '124.26'
Set-TextValue 'E49' '  +2.01%  '
Set-TextValue 'E50' '  +0.03%  '
Set-TextValue 'D51' '0.0517'
Set-TextValue 'E51' '  -1.49%  '
